# Updated working hours log: add a new entry for 23.9.2025.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 13 (row 12 intentionally left blank, same gap pattern as row 6).
$ws.Range("A13").Value = "23.9.2025"

# Match the time-of-day number format already used by the other data rows
# so the new cells get styled the same way (style index 1 / numFmtId 18).
$ws.Range("B13:E13").NumberFormat = $ws.Range("B2:E2").NumberFormat

$ws.Range("B13").Value = 0.60416666666666663
$ws.Range("C13").Value = 0.66666666666666663
$ws.Range("D13").Value = 0.72916666666666663
$ws.Range("E13").Value = 0.82291666666666663

# Move the active selection the way the author left it after editing.
$ws.Range("F16").Select()
